# Update gh-pages output data (view counts, a new event row, and a
# corrected duplicate row) across all four worksheets.

$wb = $excel.ActiveWorkbook

# --- 展览 : refreshed "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 52
$ws.Range("F5").Value = 1244
$ws.Range("F6").Value = 1708
$ws.Range("F8").Value = 564
$ws.Range("F9").Value = 2495
$ws.Range("F10").Value = 723
$ws.Range("F11").Value = 569
$ws.Range("F12").Value = 562
$ws.Range("F13").Value = 20
$ws.Range("F15").Value = 331
$ws.Range("F16").Value = 213
$ws.Range("F18").Value = 2106
$ws.Range("F21").Value = 2
$ws.Range("F22").Value = 2606
$ws.Range("F25").Value = 515
$ws.Range("F27").Value = 241
$ws.Range("F28").Value = 241
$ws.Range("F29").Value = 324
$ws.Range("F30").Value = 1767
$ws.Range("F36").Value = 4556
$ws.Range("F37").Value = 126

# --- 演出 : refreshed "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 406
$ws.Range("F5").Value = 4195
$ws.Range("F8").Value = 62
$ws.Range("F11").Value = 66
$ws.Range("F14").Value = 311
$ws.Range("F16").Value = 18
$ws.Range("F18").Value = 153
$ws.Range("F20").Value = 273
$ws.Range("F26").Value = 248
$ws.Range("F28").Value = 254
$ws.Range("F32").Value = 8
$ws.Range("F35").Value = 23

# New event row appended at the bottom (row 38)
$ws.Range("A37").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value = 37
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "2024-07-19"
$ws.Range("C38").Value = "上海·《你的名字》《天气之子》《铃芽之旅》——新海诚动漫三部曲钢琴演奏会"
$ws.Range("D38").Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws.Range("E38").Value = "2024.07.19 19:30-07.19 21:30"
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 80
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=83479"
$ws.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202403/GpyueuYA1711508106584.jpeg"

# --- 本地生活 : refreshed "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1416
$ws.Range("F5").Value = 1785
$ws.Range("F6").Value = 515
$ws.Range("F7").Value = 42
$ws.Range("F8").Value = 182

# --- 全部类型 : refreshed "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1416
$ws.Range("F6").Value = 52
$ws.Range("F9").Value = 1244
$ws.Range("F10").Value = 1708
$ws.Range("F12").Value = 62
$ws.Range("F15").Value = 564
$ws.Range("F16").Value = 2495
$ws.Range("F17").Value = 723
$ws.Range("F18").Value = 569
$ws.Range("F19").Value = 562
$ws.Range("F20").Value = 20
$ws.Range("F22").Value = 331
$ws.Range("F23").Value = 66
$ws.Range("F27").Value = 18
$ws.Range("F29").Value = 2106
$ws.Range("F31").Value = 153
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = 273
$ws.Range("F35").Value = 515
$ws.Range("F36").Value = 182
$ws.Range("F40").Value = 242
$ws.Range("F41").Value = 1767
$ws.Range("F42").Value = 248
$ws.Range("F46").Value = 4556
$ws.Range("F47").Value = 126
$ws.Range("F49").Value = 23

# Row 24 previously listed a cancelled "Walk Off The Earth" show;
# it is replaced by the Ghibli exhibition (whose own row shifts down
# into the newly freed row 25, which now lists the UP! concert).
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "2024-04-12"
$ws.Range("C24").Value = "上海·吉卜力工作室物语-沉浸式艺术展全球首站"
$ws.Range("D24").Value = "龙台路10号2F 上海国际传媒港艺术中心"
$ws.Range("E24").Value = "2024.04.12 10:00-05.12 20:00"
$ws.Range("F24").Value = 213
$ws.Range("G24").Value = 158
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=83036"
$ws.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202403/aZoum5Hd1710472525792.jpeg"

$ws.Range("C25").Value = "上海·奇迹の闪耀 「UP!」巡回动漫演唱会"
$ws.Range("D25").Value = "北京西路1700号 云峰剧院"
$ws.Range("E25").Value = "2024.04.12 19:30-04.12 21:30"
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 126
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=82427"
$ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202403/HvxHPz981709707512970.jpeg"
